# Updates cryptocurrency price/volume figures in column D (Price) and
# column E (Volume(1h)) to match the latest scrape, cell by cell.
#
# Column D values are stored as *text* in the source sheet (they use
# dotted thousands separators like '61.946.78', or are free-form text
# like currency-ish decimals). When a replacement text still looks like
# a plain decimal number (e.g. '406.40'), a leading apostrophe is used
# (the normal Excel quote-prefix convention) so the value commits as text
# instead of being auto-converted to a numeric type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.946.78'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '3.414.82'
$ws.Range('E3').Value = '  -2.87%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''406.40'
$ws.Range('E5').Value = '  -1.18%  '
$ws.Range('D6').Value = '''134.64'
$ws.Range('E6').Value = '  +3.94%  '
$ws.Range('E7').Value = '  -1.16%  '
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '''0.684'
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('E10').Value = '  -6.14%  '
$ws.Range('D11').Value = '''42.60'
$ws.Range('E11').Value = '  -1.18%  '
$ws.Range('E12').Value = '  -1.08%  '
$ws.Range('D13').Value = '''8.42'
$ws.Range('E13').Value = '  -3.61%  '
$ws.Range('D14').Value = '''19.89'
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('D15').Value = '3.434.07'
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('D16').Value = '61.954.22'
$ws.Range('E16').Value = '  -1.63%  '
$ws.Range('E17').Value = '  -3.14%  '
$ws.Range('D18').Value = '''11.01'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('E19').Value = '  -5.46%  '
$ws.Range('D20').Value = '''3.18'
$ws.Range('E20').Value = '  -5.08%  '
$ws.Range('D21').Value = '''84.79'
$ws.Range('E21').Value = '  +3.34%  '
$ws.Range('D22').Value = '''312.80'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('D23').Value = '''12.84'
$ws.Range('E23').Value = '  -1.76%  '
$ws.Range('E24').Value = '  -1.28%  '
$ws.Range('D25').Value = '''4.74'
$ws.Range('E25').Value = '  +8.54%  '
$ws.Range('D26').Value = '''29.61'
$ws.Range('E26').Value = '  -3.11%  '
$ws.Range('D27').Value = '''8.13'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').Value = '''2.82'
$ws.Range('E28').Value = '  +4.52%  '
$ws.Range('D29').Value = '''7.52'
$ws.Range('E29').Value = '  -3.42%  '
$ws.Range('D30').Value = '''0.173'
$ws.Range('E30').Value = '  -4.35%  '
$ws.Range('D31').Value = '''0.115'
$ws.Range('E31').Value = '  -3.20%  '
$ws.Range('D32').Value = '''42.35'
$ws.Range('E32').Value = '  -2.16%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Value = '''11.30'
$ws.Range('E34').Value = '  -6.65%  '
$ws.Range('D35').Value = '''0.0481'
$ws.Range('E35').Value = '  -2.56%  '
$ws.Range('D36').Value = '''51.75'
$ws.Range('E36').Value = '  -1.62%  '
$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('D38').Value = '''3.40'
$ws.Range('E38').Value = '  -4.73%  '
$ws.Range('D39').Value = '''2.93'
$ws.Range('E39').Value = '  -2.86%  '
$ws.Range('D40').Value = '''0.303'
$ws.Range('E40').Value = '  +5.34%  '
$ws.Range('D41').Value = '''137.93'
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('E42').Value = '  -1.61%  '
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').Value = '''4.03'
$ws.Range('E44').Value = '  +1.46%  '
$ws.Range('E45').Value = '  -5.84%  '
$ws.Range('E46').Value = '  -1.73%  '
$ws.Range('D47').Value = '''21.38'
$ws.Range('E47').Value = '  -4.63%  '
$ws.Range('D48').Value = '2.117.15'
$ws.Range('E48').Value = '  -4.71%  '
$ws.Range('D49').Value = '''2.30'
$ws.Range('E49').Value = '  -3.54%  '
$ws.Range('D50').Value = '''1.89'
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('D51').Value = '''0.0347'
$ws.Range('E51').Value = '  +2.79%  '
